$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the two "Descrever Totalmente..." tasks from column A (Todo)
# to the end of column D (Done), reflecting their completion.
$ws.Range("D9").Value = $ws.Range("A2").Value2
$ws.Range("D10").Value = $ws.Range("A3").Value2

$ws.Range("A2").ClearContents()
$ws.Range("A3").ClearContents()

$ws.Range("D10").Select()
